# Update the "Sciences" sheet: row 4 used to describe the "QuranRevelationOccasions"
# entry (sourced from ketabonline.com) and is replaced with a
# "QuranInterpretationErrata" entry (sourced from an app.box.com share), the
# "About" row's edition date is bumped, and the workbook is left to recalc.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sciences")

# Row 5 (ScienceMinor record #4, "QuranInterpretation"): give the generic
# "التفسير" label the specific, named reference it now stands for.
$ws.Range("H5").Value = "«التفسير الوسيط للقرآن الكريم»"

# Row 4 (ScienceMinor record #3): swap the Quran "Revelation Occasions" entry
# for an "Interpretation Errata" entry for that same tafsir (icon name and all
# four store/web links updated too).
$ws.Range("G4").Value = "QuranInterpretationErrata"
$ws.Range("H4").Value = "«التفسير الوسيط...» {البقرة:61-74}"
$ws.Range("I4").Value = "exclamation"
$ws.Range("J4").Value = "https://app.box.com/s/yl8ow6iu7y1moq0jcb01zfck3kglyra3"
$ws.Range("K4").Value = "https://app.box.com/s/yl8ow6iu7y1moq0jcb01zfck3kglyra3"
$ws.Range("L4").Value = "https://app.box.com/s/yl8ow6iu7y1moq0jcb01zfck3kglyra3"
$ws.Range("M4").Value = "https://app.box.com/s/yl8ow6iu7y1moq0jcb01zfck3kglyra3"

# Row 60 ("About" / ContentEdition record): bump the published edition date.
$ws.Range("H60").Value = "طبعة @ 2023/05/21 م - 1444/11/01 هـ"

# Recalculate so dependent formulas (B4, B70 subtotal, Complements!F2 date,
# etc.) refresh against the new values.
$excel.CalculateFull()
